# "Verify_AddMoreItems" test update: the CostCenter sample row's ItemName
# changes from "Computers" to "DeskTops", and the sheet's active selection
# moves from K7 to H8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CostCenter")

# Update the ItemName cell in the sample data row
$ws.Range("D2").Value = "DeskTops"

# Move the active selection to H8
$ws.Range("H8").Select()

# Keep the workbook window sizing in sync with the authored view state
$wb.Windows.Item(1).Width = 15435
